$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3.015327856792667
$ws.Range("F2").Value = 1.908671011575838
$ws.Range("E3").Value = 7.199849539655979
$ws.Range("F3").Value = 2.729178932624134
$ws.Range("E4").Value = 7.199849539655979
$ws.Range("F4").Value = 2.729178932624134
$ws.Range("E5").Value = 7.199849539655979
$ws.Range("F5").Value = 2.729178932624134
$ws.Range("E6").Value = 7.199849539655979
$ws.Range("F6").Value = 2.729178932624134
$ws.Range("E7").Value = 9.628352970030249
$ws.Range("F7").Value = 1.035516608633646
$ws.Range("E8").Value = 9.628352970030249
$ws.Range("F8").Value = 1.035516608633646
$ws.Range("E9").Value = 9.628352970030249
$ws.Range("F9").Value = 1.035516608633646
$ws.Range("E10").Value = 9.628352970030249
$ws.Range("F10").Value = 1.035516608633646
$ws.Range("E11").Value = 8.414101254843114
$ws.Range("F11").Value = 1.88234777062889
$ws.Range("E12").Value = 8.414101254843114
$ws.Range("F12").Value = 1.88234777062889
$ws.Range("E13").Value = 8.414101254843114
$ws.Range("F13").Value = 1.88234777062889
$ws.Range("E14").Value = 8.414101254843114
$ws.Range("F14").Value = 1.88234777062889
$ws.Range("E15").Value = 8.414101254843114
$ws.Range("F15").Value = 1.88234777062889
$ws.Range("E16").Value = 8.414101254843114
$ws.Range("F16").Value = 1.88234777062889
$ws.Range("E17").Value = 12.05685640040452
$ws.Range("F17").Value = 2.729178932624134
$ws.Range("E18").Value = 12.05685640040452
$ws.Range("F18").Value = 2.729178932624134
$ws.Range("E19").Value = 12.05685640040452
$ws.Range("F19").Value = 2.729178932624134
$ws.Range("E20").Value = 12.05685640040452
$ws.Range("F20").Value = 2.729178932624134
$ws.Range("E21").Value = 14.48535983077879
$ws.Range("F21").Value = 1.035516608633646
$ws.Range("E22").Value = 14.48535983077879
$ws.Range("F22").Value = 1.035516608633646
$ws.Range("E23").Value = 14.48535983077879
$ws.Range("F23").Value = 1.035516608633646
$ws.Range("E24").Value = 14.48535983077879
$ws.Range("F24").Value = 1.035516608633646
$ws.Range("E25").Value = 13.27110811559165
$ws.Range("F25").Value = 1.88234777062889
$ws.Range("E26").Value = 13.27110811559165
$ws.Range("F26").Value = 1.88234777062889
$ws.Range("E27").Value = 13.27110811559165
$ws.Range("F27").Value = 1.88234777062889
$ws.Range("E28").Value = 13.27110811559165
$ws.Range("F28").Value = 1.88234777062889
$ws.Range("E29").Value = 13.27110811559165
$ws.Range("F29").Value = 1.88234777062889
$ws.Range("E30").Value = 13.27110811559165
$ws.Range("F30").Value = 1.88234777062889
$ws.Range("E31").Value = 18.30929365700822
$ws.Range("F31").Value = 2.735676360711819
$ws.Range("E32").Value = 18.30929365700822
$ws.Range("F32").Value = 2.735676360711819
$ws.Range("E33").Value = 18.30929365700822
$ws.Range("F33").Value = 2.735676360711819
$ws.Range("E34").Value = 18.30929365700822
$ws.Range("F34").Value = 2.735676360711819
$ws.Range("E35").Value = 20.54373743320868
$ws.Range("F35").Value = 0.9822766714948948
$ws.Range("E36").Value = 20.54373743320868
$ws.Range("F36").Value = 0.9822766714948948
$ws.Range("E37").Value = 20.54373743320868
$ws.Range("F37").Value = 0.9822766714948948
$ws.Range("E38").Value = 20.54373743320868
$ws.Range("F38").Value = 0.9822766714948948
$ws.Range("E39").Value = 19.42651554510845
$ws.Range("F39").Value = 1.858976516103357
$ws.Range("E40").Value = 19.42651554510845
$ws.Range("F40").Value = 1.858976516103357
$ws.Range("E41").Value = 19.42651554510845
$ws.Range("F41").Value = 1.858976516103357
$ws.Range("E42").Value = 19.42651554510845
$ws.Range("F42").Value = 1.858976516103357
$ws.Range("E43").Value = 19.42651554510845
$ws.Range("F43").Value = 1.858976516103357
$ws.Range("E44").Value = 19.42651554510845
$ws.Range("F44").Value = 1.858976516103357
$ws.Range("E45").Value = 16.28143135629049
$ws.Range("F45").Value = 1.861785305612455
$ws.Range("E46").Value = 16.28143135629049
$ws.Range("F46").Value = 1.861785305612455
$ws.Range("E47").Value = 21.87779108062992
$ws.Range("F47").Value = 1.858999925318772
$ws.Range("E48").Value = 21.87779108062992
$ws.Range("F48").Value = 1.858999925318772
$ws.Range("E49").Value = 21.87779108062992
$ws.Range("F49").Value = 1.858999925318772
$ws.Range("E50").Value = 21.87779108062992
$ws.Range("F50").Value = 1.858999925318772
